$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "value" column (C2:C26) from 0.25 to 0.1
$ws.Range("C2:C26").Value = 0.1

# Update the active selection to match the target state (C13)
$ws.Range("C13").Select()
